$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows("13:13").Insert()

# --- Content fix-ups ---

# Row 10 (Objetivos): replace B/C text with the new objectives paragraph
$ws.Range("B10:C10").Value = 'Introduzir os princípios básicos da genética com apresentação minuciosa dos importantes conceitos da genética clássica, genética molecular e genômica associada às técnicas e aplicações da genética molecular.'

# Row 13 (new): professor name now stands alone on its own row
$ws.Range("B13:C13").Value = '8711290 - Elisson Antônio da Costa Romanel'

# Row 14 (Programa resumido): replace "Semestral" with full summary text
$ws.Range("B14:C14").Value = 'Introdução à genética; Estrutura e Replicação Molecular do DNA; Transcrição, Tradução e Código Genético; Mutação e Reparo do DNA; Regulação da Expressão Gênica; Genômica e Bioinformática; Técnicas e Aplicações da Genética Molecular.'

# Row 16 (Programa): replace stray date text with the full Portuguese syllabus
$ws.Range("B16:C16").Value = '1. Introdução à Genética.2. Nucleotídeos e Estrutura do DNA3. Genes e Cromossomos4. Replicação do DNA5. Transcrição e Processamento do RNA6. Código genético e Tradução7. Clonagem de DNA8. Genômica9. Bioinformática10. Regulação da Expressão Gênica11. Elementos Genéticos Transponíveis12. Mutação, Reparo de DNA e Recombinação13. Técnicas e Aplicações da Genética Molecular'

# Row 19 (Metodo): replace professor-name leftover with grading notes text
$ws.Range("B19:C19").Value = 'Notas - N distribuído no semestre. A composição das "N" fica critério do docente.'

# Row 20 (Criterio): replace grading notes text with MF formula text
$ws.Range("B20:C20").Value = 'MF = MF = (somatório de N)/número de N (adequando o valor de N, quando houver peso distinto para as Ns)'

# Row 21 (Norma de recuperacao): replace MF formula text with NF formula text
$ws.Range("B21:C21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0.'

# Row 22 (Bibliografia): replace NF formula leftover text with the real bibliography
$ws.Range("B22:C22").Value = '- Fundamentos de Genética. Peter Snustad e Michael Simmons (2013 Sexta Edição). Editora Guanabara Koogan- Principles of Genetics. Dr. Peter Snustad and Michael Simmons (2016 7th Edition). John Wiley and Sons- Genética: Um enfoque conceitual. Benjamin Pierce (2016 – Quinta Edição). Editora Guanabara Koogan.- Genetics: A Conceptual Approach. Benjamin Pierce (2016 – Sixth Edition). W.H. Freeman and Company.- Introdução à genética. Lewontin, Griffiths, Carroll e Wessler. (2013 – Décima Edição). Guanabara Koogan.- An Introduction to Genetic Analysis. Anthony Griffiths, Susan Wessler, Sean Carroll, and John Doebley (2015 Eleventh Edition). W. H. Freeman'

